$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 1: "Christ MacEachern" -> "Chris MacEachern" in the "Group members"
# paragraph, typed interactively so Word splits the single run into three:
#   "Matthew Hutchinson (Project Ma" | "nager), Nimna Ekanayaka, Chris " | "MacEachern."
# -----------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Matthew Hutchinson (Project Manager), Nimna Ekanayaka, Christ MacEachern.")
$start1 = $rng1.Start

# Bookmark the two future run-boundaries so the edit below leaves the run split
# in place (matching runs created by interactive typing) instead of being
# re-coalesced with its neighbour.
$b1a = $d.Range($start1 + 30, $start1 + 30)
$d.Bookmarks.Add("ZZ_Split1", $b1a) | Out-Null
$b1b = $d.Range($start1 + 62, $start1 + 62)
$d.Bookmarks.Add("ZZ_Split2", $b1b) | Out-Null

# Delete the "t" in "Christ " (offset 60 relative to the sentence start).
$tRange = $d.Range($start1 + 60, $start1 + 61)
$tRange.Text = ""

$d.Bookmarks("ZZ_Split1").Delete()
$d.Bookmarks("ZZ_Split2").Delete()

# -----------------------------------------------------------------------
# Edit 2: move the "_GoBack" bookmark from after "...connected to" down to
# the "See previous message(s)..." bullet, where "message" becomes
# "messages" (typed interactively, splitting that run into three as well).
# -----------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$rng2 = $d.Content
$rng2.Find.Execute("See previous message sent to users on their friends list.")
$start2 = $rng2.Start

# Type the "s" after "message".
$insPos = $d.Range($start2 + 20, $start2 + 20)
$insPos.InsertBefore("s")

# Bookmark around the freshly typed "s" so it ends up as its own run.
$preS = $d.Range($start2 + 20, $start2 + 20)
$d.Bookmarks.Add("ZZ_PreS", $preS) | Out-Null
$postS = $d.Range($start2 + 21, $start2 + 21)
$d.Bookmarks.Add("ZZ_PostS", $postS) | Out-Null
$d.Bookmarks("ZZ_PreS").Delete()
$d.Bookmarks("ZZ_PostS").Delete()

# Drop "_GoBack" back in, right after the "s" -- this is where Word leaves it
# after the most recent edit.
$bmPos = $d.Range($start2 + 21, $start2 + 21)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null

Write-Host "Edits applied"
